$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Layout tweaks ---
# Column B width: 10.5 -> 11 (ColumnWidth API adds a constant 5/6 padding vs xml width)
$ws.Columns("B").ColumnWidth = 11 - 0.8333333333333334
# Row 2 height: 28.8 -> 13.55
$ws.Rows(2).RowHeight = 13.55

# --- Apply the light (white) fill used across the whole table (header + data, rows 1-12, cols A-D) ---
$ws.Range("A1:D12").Interior.ColorIndex = 2

# --- Fill in the "Trained Model (Sprint 2)" results column (D) ---
# Make sure these are stored as text (same convention as the other "number-looking" columns)
$ws.Range("D3:D12").NumberFormat = "@"

$ws.Range("D3").Value = "1.0000100135803223"
$ws.Range("D4").Value = "1.0000089406967163"
$ws.Range("D5").Value = "1.0000100135803223"
$ws.Range("D6").Value = "1.0000035762786865"
$ws.Range("D7").Value = "0.9997366070747375"
$ws.Range("D8").Value = "0.9995201826095581"
$ws.Range("D9").Value = "0.9673290252685547"
$ws.Range("D10").Value = "0.96995609998703"
$ws.Range("D11").Value = "0.5667021870613098"
$ws.Range("D12").Value = "0.999575674533844"
